$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-parsed as a plain number
# but must remain text (matching the source data's string type), so force
# text format before writing the value.
$textCells = @("D5", "D6", "D8", "D11", "D12", "D13", "D16", "D19", "D20", "D21", "D22", "D24", "D27", "D29", "D30", "D33", "D37", "D39", "D40", "D41", "D43", "D44", "D45", "D46", "D48", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "59.519.34"
$ws.Range("E2").Value = "  -2.42%  "
$ws.Range("D3").Value = "2.302.80"
$ws.Range("E3").Value = "  -3.16%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "540.09"
$ws.Range("E5").Value = "  -1.64%  "
$ws.Range("D6").Value = "127.48"
$ws.Range("E6").Value = "  -5.37%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "0.567"
$ws.Range("E8").Value = "  -3.99%  "
$ws.Range("D9").Value = "2.300.94"
$ws.Range("E9").Value = "  -3.11%  "
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("D11").Value = "5.51"
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("D12").Value = "0.150"
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("D13").Value = "0.330"
$ws.Range("E13").Value = "  -3.36%  "
$ws.Range("D14").Value = "2.714.69"
$ws.Range("E14").Value = "  -3.39%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "59.526.18"
$ws.Range("E15").Value = "  -2.29%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "23.02"
$ws.Range("E16").Value = "  -5.26%  "
$ws.Range("E17").Value = "  -3.03%  "
$ws.Range("D18").Value = "2.300.30"
$ws.Range("E18").Value = "  -3.18%  "
$ws.Range("D19").Value = "10.38"
$ws.Range("E19").Value = "  -4.45%  "
$ws.Range("D20").Value = "4.01"
$ws.Range("E20").Value = "  -5.80%  "
$ws.Range("D21").Value = "309.28"
$ws.Range("E21").Value = "  -3.32%  "
$ws.Range("D22").Value = "6.47"
$ws.Range("E22").Value = "  -6.64%  "
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("D24").Value = "62.95"
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("E25").Value = "  -3.49%  "
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").Value = "7.67"
$ws.Range("E27").Value = "  -6.92%  "
$ws.Range("E28").Value = "  -2.31%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "171.42"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("B30").Value = "SuiNetwork"
$ws.Range("C30").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D30").Value = "1.18"
$ws.Range("E30").Value = "  +2.46%  "
$ws.Range("E31").Value = "  -3.68%  "
$ws.Range("E32").Value = "  -5.83%  "
$ws.Range("D33").Value = "5.77"
$ws.Range("E33").Value = "  -3.53%  "
$ws.Range("E34").Value = "  -3.52%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  -7.81%  "
$ws.Range("D37").Value = "17.65"
$ws.Range("E37").Value = "  -2.77%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "3.96"
$ws.Range("E39").Value = "  -6.43%  "
$ws.Range("D40").Value = "310.99"
$ws.Range("E40").Value = "  -5.54%  "
$ws.Range("D41").Value = "37.56"
$ws.Range("E41").Value = "  -2.47%  "
$ws.Range("E42").Value = "  -5.86%  "
$ws.Range("D43").Value = "135.65"
$ws.Range("E43").Value = "  -7.56%  "
$ws.Range("D44").Value = "3.39"
$ws.Range("E44").Value = "  -3.48%  "
$ws.Range("D45").Value = "0.0934"
$ws.Range("E45").Value = "  -2.69%  "
$ws.Range("D46").Value = "0.568"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D48").Value = "18.42"
$ws.Range("E48").Value = "  -6.49%  "
$ws.Range("E49").Value = "  +22.89%  "
$ws.Range("E50").Value = "  -2.27%  "
$ws.Range("D51").Value = "10.98"
$ws.Range("E51").Value = "  -0.40%  "
